$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header (row 1), shifting the existing
# data rows (2..108) down to (3..109). Doing this by copying each row's
# range into the row below it (bottom-up) rather than using Rows.Insert
# keeps the existing style table untouched (no new cell formats appear).
$lastRow = 108
for ($r = $lastRow; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $src.Copy($dst)
}
$excel.CutCopyMode = $false

# Fill the freed-up row 2 with the new "All" / "2 stage" entry. Give A2
# the same bold header formatting as A1 (copy its format only), while B2
# keeps the plain formatting it inherited from the old row 2 above.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 1).Value = "All"
$ws.Cells.Item(2, 2).Value = "2 stage"

# Move the selection to B3, matching the resulting workbook state.
$ws.Range("B3").Select()
